$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) is treated as text so values like "1.00" or "0.150"
# keep their exact textual representation instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "54.512.92"
$ws.Range("E2").Value = "  +1.04%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "2.284.34"
$ws.Range("E3").Value = "  +1.53%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "503.47"
$ws.Range("E5").Value = "  +2.36%  "

$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").Value = "130.13"
$ws.Range("E6").Value = "  +2.50%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "0.528"
$ws.Range("E8").Value = "  +0.65%  "

$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "2.292.22"
$ws.Range("E9").Value = "  +1.27%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.0959"
$ws.Range("E10").Value = "  +1.07%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.153"
$ws.Range("E11").Value = "  +1.07%  "

$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").Value = "0.337"
$ws.Range("E12").Value = "  +4.21%  "

$ws.Range("B13").Value = "Toncoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D13").Value = "4.92"
$ws.Range("E13").Value = "  +6.17%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "23.12"
$ws.Range("E14").Value = "  +6.61%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.690.34"
$ws.Range("E15").Value = "  +1.53%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "54.520.08"
$ws.Range("E16").Value = "  +1.23%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0000131"
$ws.Range("E17").Value = "  +1.91%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.289.53"
$ws.Range("E18").Value = "  +2.80%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "10.32"
$ws.Range("E19").Value = "  +3.35%  "

$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "4.16"
$ws.Range("E20").Value = "  +2.50%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "305.95"
$ws.Range("E21").Value = "  +2.10%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "6.32"
$ws.Range("E22").Value = "  -1.13%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "60.57"
$ws.Range("E24").Value = "  -2.08%  "

$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D25").Value = "0.995"
$ws.Range("E25").Value = "  -1.91%  "

$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "0.150"
$ws.Range("E26").Value = "  +1.86%  "

$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "7.39"
$ws.Range("E27").Value = "  +5.07%  "

$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "172.36"
$ws.Range("E28").Value = "  +3.85%  "

$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0708"
$ws.Range("E29").Value = "  +4.89%  "

$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "6.05"
$ws.Range("E30").Value = "  +3.95%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.62"
$ws.Range("E31").Value = "  +1.85%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.10"
$ws.Range("E32").Value = "  +3.66%  "

$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "17.96"
$ws.Range("E34").Value = "  +2.03%  "

$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "0.994"
$ws.Range("E35").Value = "  -0.18%  "

$ws.Range("B36").Value = "SuiNetwork"
$ws.Range("C36").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D36").Value = "0.937"
$ws.Range("E36").Value = "  +6.11%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "1.20"
$ws.Range("E37").Value = "  +2.58%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "3.76"
$ws.Range("E38").Value = "  +3.37%  "

$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "36.14"
$ws.Range("E39").Value = "  +1.31%  "

$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "0.376"
$ws.Range("E40").Value = "  +1.74%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "1.41"
$ws.Range("E41").Value = "  +2.14%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "5.07"
$ws.Range("E42").Value = "  +3.30%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "3.41"
$ws.Range("E43").Value = "  +2.21%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "125.42"
$ws.Range("E44").Value = "  +0.61%  "

$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "0.0493"
$ws.Range("E45").Value = "  +2.81%  "

$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "246.99"
$ws.Range("E46").Value = "  +5.38%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "0.0897"
$ws.Range("E47").Value = "  +1.40%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.547"
$ws.Range("E48").Value = "  +1.90%  "

$ws.Range("B49").Value = "Polygon"
$ws.Range("C49").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D49").Value = "0.375"
$ws.Range("E49").Value = "  +1.67%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0207"
$ws.Range("E50").Value = "  +3.22%  "

$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "10.82"
$ws.Range("E51").Value = "  +0.78%  "

